# issue #16 Minimapa a velka mapa
# Append two new paragraphs after the last paragraph of the document
# ("Anotace [Header] ..."), right before the section break, describing
# Culling mask and RenderTexture (used for the minimap).

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$openXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = @'
<w:p><w:r><w:t xml:space="preserve">Culling mask v kameře </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> určení co kamera vlastně může vidět</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">RenderTexture v kameře </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> soubor, do kterého se ukládá, co zrovna kamera vidí </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> využití například v UI komponentě raw image (oproti klasickému Image nepoužívá sprite, ale texturu) = použítí k minimapě</w:t></w:r></w:p>
'@

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $openXmlNs + '><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$insertionPoint.InsertXML($packageXml)
